$wb = $excel.ActiveWorkbook

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23534.746
$ws.Range("I32").Value = 25680.268
$ws.Range("K32").Value = 25680.268
$ws.Range("M32").Value = -25393.268
$ws.Range("H61").Value = 5367.1064
$ws.Range("I61").Value = 3582.4595
$ws.Range("J61").Value = 11970.3
$ws.Range("K61").Value = 3582.4595
$ws.Range("L61").Value = 11970.3
$ws.Range("M61").Value = -3370.4595
$ws.Range("N61").Value = -12394.3
$ws.Range("H74").Value = 6054.089
$ws.Range("I74").Value = 4817.6562
$ws.Range("J74").Value = 9097.615
$ws.Range("K74").Value = 4817.6562
$ws.Range("L74").Value = 9097.615
$ws.Range("M74").Value = -3943.6562
$ws.Range("N74").Value = -10845.615
$ws.Range("H77").Value = 6054.089
$ws.Range("I77").Value = 4817.6562
$ws.Range("J77").Value = 9097.615
$ws.Range("K77").Value = 24088.281
$ws.Range("L77").Value = 45488.075
$ws.Range("M77").Value = -19720.281
$ws.Range("N77").Value = -54224.075
$ws.Range("H132").Value = 5311.972
$ws.Range("I132").Value = 1782.2609
$ws.Range("J132").Value = 11556.846
$ws.Range("K132").Value = 5346.7827
$ws.Range("L132").Value = 34670.538
$ws.Range("M132").Value = -2816.7827
$ws.Range("N132").Value = -39730.538
$ws.Range("H136").Value = 5367.1064
$ws.Range("I136").Value = 3582.4595
$ws.Range("J136").Value = 11970.3
$ws.Range("K136").Value = 10747.3785
$ws.Range("L136").Value = 35910.89999999999
$ws.Range("M136").Value = -8197.378499999999
$ws.Range("N136").Value = -41010.89999999999

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 756.6799999999999
$ws.Range("I20").Value = 768.4761999999999
$ws.Range("J20").Value = 694.75
$ws.Range("K20").Value = 768.4761999999999
$ws.Range("L20").Value = 694.75
$ws.Range("M20").Value = -521.4761999999999
$ws.Range("N20").Value = -1188.75
$ws.Range("H134").Value = 3396.5789
$ws.Range("I134").Value = 3594.5833
$ws.Range("J134").Value = 3057.1428
$ws.Range("K134").Value = 10783.7499
$ws.Range("L134").Value = 9171.428400000001
$ws.Range("M134").Value = -8248.749899999999
$ws.Range("N134").Value = -14241.4284

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1663.3846
$ws.Range("I16").Value = 1177.75
$ws.Range("J16").Value = 1879.2222
$ws.Range("K16").Value = 1177.75
$ws.Range("L16").Value = 1879.2222
$ws.Range("M16").Value = -890.75
$ws.Range("N16").Value = -2453.2222
$ws.Range("H31").Value = 3809.3394
$ws.Range("I31").Value = 4415.5938
$ws.Range("J31").Value = 3001
$ws.Range("K31").Value = 4415.5938
$ws.Range("L31").Value = 3001
$ws.Range("M31").Value = -4120.5938
$ws.Range("N31").Value = -3591
$ws.Range("H34").Value = 3809.3394
$ws.Range("I34").Value = 4415.5938
$ws.Range("J34").Value = 3001
$ws.Range("K34").Value = 4415.5938
$ws.Range("L34").Value = 3001
$ws.Range("M34").Value = -4213.5938
$ws.Range("N34").Value = -3405
$ws.Range("H58").Value = 1820161.1
$ws.Range("I58").Value = 2598717.5
$ws.Range("J58").Value = 3529.8667
$ws.Range("K58").Value = 2598717.5
$ws.Range("L58").Value = 3529.8667
$ws.Range("M58").Value = -2598514.5
$ws.Range("N58").Value = -3935.8667
$ws.Range("H113").Value = 1663.3846
$ws.Range("I113").Value = 1177.75
$ws.Range("J113").Value = 1879.2222
$ws.Range("K113").Value = 1177.75
$ws.Range("L113").Value = 1879.2222
$ws.Range("M113").Value = 992.25
$ws.Range("N113").Value = -6219.2222
$ws.Range("H132").Value = 1926
$ws.Range("I132").Value = 1603.6842
$ws.Range("J132").Value = 2538.4
$ws.Range("K132").Value = 4811.0526
$ws.Range("L132").Value = 7615.200000000001
$ws.Range("M132").Value = -2281.0526
$ws.Range("N132").Value = -12675.2
$ws.Range("H134").Value = 2729.423
$ws.Range("I134").Value = 1472.3125
$ws.Range("J134").Value = 4740.8
$ws.Range("K134").Value = 4416.9375
$ws.Range("L134").Value = 14222.4
$ws.Range("M134").Value = -1881.9375
$ws.Range("N134").Value = -19292.4
$ws.Range("H136").Value = 1820161.1
$ws.Range("I136").Value = 2598717.5
$ws.Range("J136").Value = 3529.8667
$ws.Range("K136").Value = 7796152.5
$ws.Range("L136").Value = 10589.6001
$ws.Range("M136").Value = -7793602.5
$ws.Range("N136").Value = -15689.6001

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 39888.4
$ws.Range("I131").Value = 1719.0769
$ws.Range("K131").Value = 5157.2307
$ws.Range("M131").Value = -117.2307000000001

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6477.385
$ws.Range("I122").Value = 7555.3335
$ws.Range("J122").Value = 4052
$ws.Range("K122").Value = 22666.0005
$ws.Range("L122").Value = 12156
$ws.Range("M122").Value = -20216.0005
$ws.Range("N122").Value = -17056
$ws.Range("H132").Value = 3176.7058
$ws.Range("I132").Value = 2933.7334
$ws.Range("K132").Value = 8801.200199999999
$ws.Range("M132").Value = -6271.200199999999

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2350
$ws.Range("I82").Value = 2160
$ws.Range("J82").Value = 2666.6667
$ws.Range("K82").Value = 2160
$ws.Range("L82").Value = 2666.6667
$ws.Range("M82").Value = -1799
$ws.Range("N82").Value = -3388.6667
$ws.Range("H85").Value = 2350
$ws.Range("I85").Value = 2160
$ws.Range("J85").Value = 2666.6667
$ws.Range("K85").Value = 2160
$ws.Range("L85").Value = 2666.6667
$ws.Range("M85").Value = -912
$ws.Range("N85").Value = -5162.6667
$ws.Range("H95").Value = 39900
$ws.Range("J95").Value = 39900
$ws.Range("L95").Value = 39900
$ws.Range("N95").Value = -45392
$ws.Range("H136").Value = 5252.0464
$ws.Range("I136").Value = 3571.9473
$ws.Range("J136").Value = 6582.125
$ws.Range("K136").Value = 10715.8419
$ws.Range("L136").Value = 19746.375
$ws.Range("M136").Value = -8165.841899999999
$ws.Range("N136").Value = -24846.375

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("H132").Value = 2923.5417
$ws.Range("I132").Value = 1598.5
$ws.Range("J132").Value = 4248.5835
$ws.Range("K132").Value = 4795.5
$ws.Range("L132").Value = 12745.7505
$ws.Range("M132").Value = -2265.5
$ws.Range("N132").Value = -17805.7505
$ws.Range("H136").Value = 8249.385
$ws.Range("I136").Value = 6732.1816
$ws.Range("J136").Value = 10212.823
$ws.Range("K136").Value = 20196.5448
$ws.Range("L136").Value = 30638.469
$ws.Range("M136").Value = -17646.5448
$ws.Range("N136").Value = -35738.469
$ws.Range("N97").ClearContents()
